$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the summary/header figures ---------------------------------
# Valor Mora total
$ws.Range("E11").Value = 113601
# Cant. Trabajadores
$ws.Range("C13").Value = 3
# Cant. Periodos
$ws.Range("F13").Value = 3

# --- 2) Make room for the two new detail rows ------------------------------
# Row 17 currently holds the last detail line (with the heavier "closing"
# border). Insert two blank rows right after it; this pushes the old
# signature block (rows 22-23) down to rows 24-25, matching the new
# dimension B2:J25.
$ws.Rows("18:19").Insert()

# Preserve the old "last row" look (heavier border) by copying it from the
# row that used to be last (still row 17 at this point) onto the new last
# row (19).
$ws.Range("B17:J17").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

# Row 17 is no longer the last row of the table, so it now gets the regular
# interior-row formatting (same as row 16). The brand-new row 18 gets the
# same regular formatting too.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J18").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- 3) Fill in the new worker (row 18) and the updated totals row (19) ---
# Row 18: new worker JAIR FERNANDO PEÑA, period 2508
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "71755863"
$ws.Range("D18").Value = "JAIR FERNANDO PEÑA"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 5395
$ws.Range("G18").Value = 4046384

# Row 19: LINDA CAROLINA POSADA SALINAS now also owes period 2508
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1002245132"
$ws.Range("D19").Value = "LINDA CAROLINA POSADA SALINAS"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500
